# PNAD 2009 - "correção nos dados e inicio da analise PNAD 2009"
#
# The original sheet had three label-only rows interleaved with the data
# (a "situação do domicílio" sub-header, a "grandes regiões e unidades da
# federação" sub-header, and a trailing "fonte: ..." footnote row) plus a
# mislabeled header cell ("unnamed: 1_level_1"). The fix removes those three
# rows entirely (shifting the remaining data rows up so every region lines
# up with its correct figures) and relabels the "total" column header.
#
# Row deletions are performed bottom-to-top so earlier row numbers stay valid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing footnote row ("fonte: ibge, ...")
$ws.Rows("41:41").Delete()

# Remove the "grandes regiões e unidades da federação" sub-header row
$ws.Rows("8:8").Delete()

# Remove the "situação do domicílio" sub-header row
$ws.Rows("5:5").Delete()

# Fix the column header that pandas exported as "unnamed: 1_level_1"
$ws.Range("B2").Value = "total"
